# Generate Report for Handback
# Updates the timestamp strings recorded during handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview.Range("G2").Value = "2016-12-15 04:52:28"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (L2)
$wsZhCn.Range("H2").Value = "2016-12-15 04:52:15"
$wsZhCn.Range("L2").Value = "2016-12-15 04:53:08"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (L2)
$wsDeDe.Range("H2").Value = "2016-12-15 04:52:28"
$wsDeDe.Range("L2").Value = "2016-12-15 04:53:26"
